# -----------------------------------------------------------------------------
# Scheduled market-data refresh: updates currentAveragePrice / LevePrice /
# LeveProfit columns (H-N) on the per-job Leve profit sheets with the latest
# Universalis market snapshot. Pure data overwrite, no formula/layout changes.
# -----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 42200
$ws.Range("J3").Value = 42200
$ws.Range("L3").Value = 42200
$ws.Range("N3").Value = -42428

$ws.Range("H6").Value = 1572.1111
$ws.Range("I6").Value = 1352
$ws.Range("J6").Value = 3333
$ws.Range("K6").Value = 4056
$ws.Range("L6").Value = 9999
$ws.Range("M6").Value = -3944
$ws.Range("N6").Value = -10223

$ws.Range("H9").Value = 159.08696
$ws.Range("I9").Value = 155.2381
$ws.Range("J9").Value = 199.5
$ws.Range("K9").Value = 155.2381
$ws.Range("L9").Value = 199.5
$ws.Range("M9").Value = 13.7619
$ws.Range("N9").Value = -537.5

$ws.Range("H12").Value = 113.333336
$ws.Range("I12").Value = 127.5
$ws.Range("J12").Value = 85
$ws.Range("K12").Value = 127.5
$ws.Range("L12").Value = 85
$ws.Range("M12").Value = 42.5
$ws.Range("N12").Value = -425

$ws.Range("H29").Value = 4028.6667
$ws.Range("I29").Value = 4028.6667
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 12086.0001
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -11805.0001
$ws.Range("N29").Value = $null

$ws.Range("H38").Value = 448.14285
$ws.Range("I38").Value = 106.166664
$ws.Range("J38").Value = 2500
$ws.Range("K38").Value = 318.499992
$ws.Range("L38").Value = 7500
$ws.Range("M38").Value = 53.50000799999998
$ws.Range("N38").Value = -8244

$ws.Range("H43").Value = 2923.077
$ws.Range("I43").Value = 3300
$ws.Range("J43").Value = 1666.6666
$ws.Range("K43").Value = 3300
$ws.Range("L43").Value = 1666.6666
$ws.Range("M43").Value = -3231
$ws.Range("N43").Value = -1804.6666

$ws.Range("H86").Value = 2084.913
$ws.Range("I86").Value = 2241.2632
$ws.Range("J86").Value = 1342.25
$ws.Range("K86").Value = 2241.2632
$ws.Range("L86").Value = 1342.25
$ws.Range("M86").Value = -1118.2632
$ws.Range("N86").Value = -3588.25

$ws.Range("H87").Value = 49350
$ws.Range("J87").Value = 49350
$ws.Range("L87").Value = 49350
$ws.Range("N87").Value = -51846

$ws.Range("H89").Value = 2084.913
$ws.Range("I89").Value = 2241.2632
$ws.Range("J89").Value = 1342.25
$ws.Range("K89").Value = 11206.316
$ws.Range("L89").Value = 6711.25
$ws.Range("M89").Value = -5590.315999999999
$ws.Range("N89").Value = -17943.25

$ws.Range("H90").Value = 49350
$ws.Range("J90").Value = 49350
$ws.Range("L90").Value = 148050
$ws.Range("N90").Value = -160530

$ws.Range("H95").Value = 36613
$ws.Range("J95").Value = 36613
$ws.Range("L95").Value = 36613
$ws.Range("N95").Value = -42105

$ws.Range("H102").Value = 42200
$ws.Range("J102").Value = 42200
$ws.Range("L102").Value = 42200
$ws.Range("N102").Value = -48690

$ws.Range("H123").Value = 36773.6
$ws.Range("J123").Value = 36773.6
$ws.Range("L123").Value = 36773.6
$ws.Range("N123").Value = -46573.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1901.4773
$ws.Range("I61").Value = 1460.1936
$ws.Range("J61").Value = 2953.7693
$ws.Range("K61").Value = 1460.1936
$ws.Range("L61").Value = 2953.7693
$ws.Range("M61").Value = -1248.1936
$ws.Range("N61").Value = -3377.7693

$ws.Range("H74").Value = 1879.0238
$ws.Range("I74").Value = 1420
$ws.Range("J74").Value = 5275.8
$ws.Range("K74").Value = 1420
$ws.Range("L74").Value = 5275.8
$ws.Range("M74").Value = -546
$ws.Range("N74").Value = -7023.8

$ws.Range("H77").Value = 1879.0238
$ws.Range("I77").Value = 1420
$ws.Range("J77").Value = 5275.8
$ws.Range("K77").Value = 7100
$ws.Range("L77").Value = 26379
$ws.Range("M77").Value = -2732
$ws.Range("N77").Value = -35115

$ws.Range("H96").Value = 32245
$ws.Range("J96").Value = 32245
$ws.Range("L96").Value = 32245
$ws.Range("N96").Value = -37737

$ws.Range("H104").Value = 29434.2
$ws.Range("J104").Value = 29434.2
$ws.Range("L104").Value = 29434.2
$ws.Range("N104").Value = -36422.2

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = $null

$ws.Range("H136").Value = 1901.4773
$ws.Range("I136").Value = 1460.1936
$ws.Range("J136").Value = 2953.7693
$ws.Range("K136").Value = 4380.5808
$ws.Range("L136").Value = 8861.3079
$ws.Range("M136").Value = -1830.5808
$ws.Range("N136").Value = -13961.3079

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16859.143
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 16859.143
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 16859.143
$ws.Range("M31").Value = $null
$ws.Range("N31").Value = -17449.143

$ws.Range("H34").Value = 16859.143
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 16859.143
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 16859.143
$ws.Range("M34").Value = $null
$ws.Range("N34").Value = -17263.143

$ws.Range("H96").Value = 78248
$ws.Range("J96").Value = 78248
$ws.Range("L96").Value = 78248
$ws.Range("N96").Value = -83740

$ws.Range("H106").Value = 34893
$ws.Range("J106").Value = 42491.25
$ws.Range("L106").Value = 42491.25
$ws.Range("N106").Value = -45015.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 41670176
$ws.Range("I136").Value = 100002420
$ws.Range("J136").Value = 4295.143
$ws.Range("K136").Value = 300007260
$ws.Range("L136").Value = 12885.429
$ws.Range("M136").Value = -300002160
$ws.Range("N136").Value = -23085.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 266916.1
$ws.Range("J80").Value = 3050.5
$ws.Range("L80").Value = 3050.5
$ws.Range("N80").Value = -5046.5

$ws.Range("H83").Value = 266916.1
$ws.Range("J83").Value = 3050.5
$ws.Range("L83").Value = 15252.5
$ws.Range("N83").Value = -25236.5

$ws.Range("H98").Value = 46635
$ws.Range("J98").Value = 46635
$ws.Range("L98").Value = 46635
$ws.Range("N98").Value = -52625

$ws.Range("H104").Value = 32994
$ws.Range("J104").Value = 32994
$ws.Range("L104").Value = 32994
$ws.Range("N104").Value = -39982

$ws.Range("H122").Value = 1749.8422
$ws.Range("I122").Value = 1834.5
$ws.Range("J122").Value = 1710.7693
$ws.Range("K122").Value = 5503.5
$ws.Range("L122").Value = 5132.3079
$ws.Range("M122").Value = -3053.5
$ws.Range("N122").Value = -10032.3079

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 44300
$ws.Range("J106").Value = 44300
$ws.Range("L106").Value = 44300
$ws.Range("N106").Value = -46824

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 23416.666
$ws.Range("J94").Value = 23416.666
$ws.Range("L94").Value = 23416.666
$ws.Range("N94").Value = -25218.666

$ws.Range("H95").Value = 38781.332
$ws.Range("J95").Value = 38781.332
$ws.Range("L95").Value = 38781.332
$ws.Range("N95").Value = -44273.332

$ws.Range("H98").Value = 46589
$ws.Range("J98").Value = 46589
$ws.Range("L98").Value = 46589
$ws.Range("N98").Value = -52579

$ws.Range("H104").Value = 43588
$ws.Range("J104").Value = 43588
$ws.Range("L104").Value = 43588
$ws.Range("N104").Value = -50576

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = $null
